$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 153, pushing existing rows 153:239 down to 154:240.
$ws.Rows.Item(153).Insert()

# Populate the newly inserted row 153 with the new record's data.
$ws.Range("A153").Value = 5
$ws.Range("B153").Value = "Macroferia Regional de Talca"
$ws.Range("C153").Value = "Maule"
$ws.Range("D153").Value = 44529
$ws.Range("E153").Value = 7
$ws.Range("F153").Value = 100112032
$ws.Range("G153").Value = "Zapallo italiano"
$ws.Range("H153").Value = "Sin especificar"
$ws.Range("I153").Value = "Primera"
$ws.Range("J153").Value = 500
$ws.Range("K153").Value = 5000
$ws.Range("L153").Value = 5000
$ws.Range("M153").Value = 5000
$ws.Range("N153").Value = "$/caja 60 unidades"
$ws.Range("O153").Value = "Región del Maule"
$ws.Range("P153").Value = 83
$ws.Range("Q153").Value = 60
$ws.Range("R153").Value = "Hortaliza"
